$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, newD (or $null if unchanged), newE, dNeedsTextFormat
$updates = @(
    ,@(2, "36.584.94", "  +0.41%  ", $false)
    ,@(3, "1.960.38", "  +0.89%  ", $false)
    ,@(4, $null, "  +0.05%  ", $false)
    ,@(5, "244.10", "  +0.36%  ", $true)
    ,@(6, "0.617", "  +0.58%  ", $true)
    ,@(7, "60.75", "  +5.74%  ", $true)
    ,@(8, "0.999", "  -0.02%  ", $true)
    ,@(9, "0.375", "  +4.04%  ", $true)
    ,@(10, "0.0789", "  -6.87%  ", $true)
    ,@(11, $null, "  +0.31%  ", $false)
    ,@(12, "14.28", "  +5.75%  ", $true)
    ,@(13, "21.89", "  +2.18%  ", $true)
    ,@(14, "0.831", "  +2.33%  ", $true)
    ,@(15, "2.248.35", "  +0.89%  ", $false)
    ,@(16, "5.27", "  +2.08%  ", $true)
    ,@(17, "1.963.77", "  +0.60%  ", $false)
    ,@(18, "36.503.47", "  +0.33%  ", $false)
    ,@(19, "69.76", "  +0.68%  ", $true)
    ,@(20, "0.0₃0852", "  -1.34%  ", $false)
    ,@(21, "229.85", "  +0.48%  ", $true)
    ,@(22, "5.07", "  +1.33%  ", $true)
    ,@(23, $null, "  +0.00%  ", $false)
    ,@(24, "2.44", "  +3.46%  ", $true)
    ,@(25, $null, "  +2.39%  ", $false)
    ,@(26, "0.142", "  +6.72%  ", $true)
    ,@(27, "9.21", "  -0.04%  ", $true)
    ,@(28, "161.04", "  -0.61%  ", $true)
    ,@(29, "19.37", "  +0.83%  ", $true)
    ,@(30, "1.30", "  +18.71%  ", $true)
    ,@(31, "0.119", "  +1.68%  ", $true)
    ,@(32, "4.77", "  +4.04%  ", $true)
    ,@(33, "0.0615", "  -0.45%  ", $true)
    ,@(34, "4.42", "  +5.63%  ", $true)
    ,@(35, "3.52", "  +11.23%  ", $true)
    ,@(36, "2.27", "  +4.49%  ", $true)
    ,@(37, $null, "  -0.08%  ", $false)
    ,@(38, $null, "  -1.03%  ", $false)
    ,@(39, "5.50", "  -11.98%  ", $true)
    ,@(40, "0.0979", "  -1.24%  ", $true)
    ,@(42, $null, "  +1.38%  ", $false)
    ,@(43, "0.0210", "  +0.57%  ", $true)
    ,@(44, "15.93", "  -0.09%  ", $true)
    ,@(45, "1.367.78", "  +1.88%  ", $false)
    ,@(46, "88.76", "  +2.34%  ", $true)
    ,@(47, $null, "  +0.26%  ", $false)
    ,@(48, "7.15", "  -0.31%  ", $true)
    ,@(49, $null, "  +0.39%  ", $false)
    ,@(50, "45.40", "  +4.79%  ", $true)
    ,@(51, "2.138.31", "  +0.89%  ", $false)
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    $needsText = $u[3]
    if ($null -ne $dVal) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($needsText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $dVal
    }
    $ws.Cells.Item($row, 5).Value = $eVal
}

Write-Host "Done"
